$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the IPO entries for 큐로셀 (row 18), 메가터치 (row 22),
# 컨텍 (row 23) and 비아이매트릭스 (row 24).
# Delete from the bottom up so the remaining row numbers stay valid.
$rowsToDelete = @(24, 23, 22, 18)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
